$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "Influenza"
$ws.Range("E4").Value = "Influenza"
$ws.Range("E5").Value = "Influenza"
$ws.Range("E6").Value = "Influenza"
